# Add a "premium_experience" column (AR) to the Zakopane POI sheet and
# flag the one premium POI (KULIGI, row 32) as TRUE; style the header
# row with a bold font, thin box border and centered/top alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 36
$newCol = 44          # column AR
$newColLetter = "AR"

# --- 1. New column header -------------------------------------------------
$headerCell = $ws.Cells.Item(1, $newCol)
$headerCell.Value = "premium_experience"

# --- 2. premium_experience values for existing data rows ------------------
# Only row 32 (KULIGI w Zakopanem) is a premium experience.
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $newCol)
    if ($row -eq 32) {
        $cell.Value = $true
    } else {
        $cell.Value = $false
    }
}

# --- 3. Header-row styling (bold, thin box border, center/top align) ------
$headerRange = $ws.Range("A1:" + $newColLetter + "1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin box border)

Write-Host "premium_experience column added; header styled"
